# Weekly price update: insert the newest week's data ("Apio" - Celery) at the
# top of the historical data table (rows 135-136), shifting all existing
# history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 135 (pushes old rows 135:258 down to 137:260,
# automatically preserving all of their existing values/formats).
$ws.Range("A135:A136").EntireRow.Insert()

# --- New row 135: "Primera" quality ---
$ws.Range("A135").Value() = 8
$ws.Range("B135").Value() = "Terminal La Palmera de La Serena"
$ws.Range("C135").Value() = "Coquimbo"
$ws.Range("D135").Value() = 44512
$ws.Range("E135").Value() = 4
$ws.Range("F135").Value() = 100112017
$ws.Range("G135").Value() = "Apio"
$ws.Range("H135").Value() = "Americana (o)"
$ws.Range("I135").Value() = "Primera"
$ws.Range("J135").Value() = 2600
$ws.Range("K135").Value() = 7000
$ws.Range("L135").Value() = 8000
$ws.Range("M135").Value() = 7500
$ws.Range("N135").Value() = "$/docena de matas"
$ws.Range("O135").Value() = "Provincia del Elquí"
$ws.Range("P135").Value() = 1250
$ws.Range("Q135").Value() = 6
$ws.Range("R135").Value() = "Hortaliza"

# --- New row 136: "Segunda" quality ---
$ws.Range("A136").Value() = 8
$ws.Range("B136").Value() = "Terminal La Palmera de La Serena"
$ws.Range("C136").Value() = "Coquimbo"
$ws.Range("D136").Value() = 44512
$ws.Range("E136").Value() = 4
$ws.Range("F136").Value() = 100112017
$ws.Range("G136").Value() = "Apio"
$ws.Range("H136").Value() = "Americana (o)"
$ws.Range("I136").Value() = "Segunda"
$ws.Range("J136").Value() = 1300
$ws.Range("K136").Value() = 5500
$ws.Range("L136").Value() = 6000
$ws.Range("M136").Value() = 5750
$ws.Range("N136").Value() = "$/docena de matas"
$ws.Range("O136").Value() = "Provincia del Elquí"
$ws.Range("P136").Value() = 958
$ws.Range("Q136").Value() = 6
$ws.Range("R136").Value() = "Hortaliza"
